$wb = $excel.ActiveWorkbook

# "Bids" sheet: clear out all existing bid rows, keeping only the header row.
$bids = $wb.Worksheets.Item("Bids")
$bids.Range("A2:C5").ClearContents()

# "Participants" sheet: update the Base Price / Budget for Hiren and Anurag to 100.
$participants = $wb.Worksheets.Item("Participants")
$participants.Range("B3").Value = 100
$participants.Range("B4").Value = 100
